# Research/Tests/DownscalingNew.xlsx - "Enemies: Fixed some missing data and
# added default values for stats."
#
# 1. "Graph Gear Attributes": fill in the previously-missing target level
#    (A71) for the weapon-power downscaling test row; the dependent
#    formulas in B71:D71 (and the chart caches that read this range)
#    recalculate automatically.
# 2. Add a brand-new "Weapon" worksheet at the end of the workbook with
#    the default per-rarity weapon stat tables.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) Graph Gear Attributes: supply the missing target level at row 71.
# ---------------------------------------------------------------------
$wsGear = $wb.Worksheets.Item("Graph Gear Attributes")
$wsGear.Activate()
$wsGear.Range("A71").Value = 1966
$wsGear.Range("A72").Select()

# ---------------------------------------------------------------------
# 2) Add the "Weapon" worksheet (new last tab, becomes the active sheet).
# ---------------------------------------------------------------------
$wsWeapon = $wb.Worksheets.Add($null, $wb.Sheets($wb.Sheets.Count))
$wsWeapon.Name = "Weapon"

# --- Left block (Power) ---------------------------------------------
$wsWeapon.Range("A1").Value = "Level"
$wsWeapon.Range("B1").Value = 80

$wsWeapon.Range("A2").Value = "Basic"
$wsWeapon.Range("B2").Value = 551
$wsWeapon.Range("C2").Value = 609
$wsWeapon.Range("D2").Formula = "=(B2+C2)/2"
$wsWeapon.Range("E2").Formula = "=D5/D2"

$wsWeapon.Range("A3").Value = "Fine"
$wsWeapon.Range("B3").Value = 690
$wsWeapon.Range("C3").Value = 762
$wsWeapon.Range("D3").Formula = "=(B3+C3)/2"

$wsWeapon.Range("A4").Value = "Masterwork"
$wsWeapon.Range("B4").Value = 745
$wsWeapon.Range("C4").Value = 823
$wsWeapon.Range("D4").Formula = "=(B4+C4)/2"
$wsWeapon.Range("E4").Formula = "=D5/D4"

$wsWeapon.Range("A5").Value = "Rare"
$wsWeapon.Range("B5").Value = 802
$wsWeapon.Range("C5").Value = 886
$wsWeapon.Range("D5").Formula = "=(B5+C5)/2"
$wsWeapon.Range("E5").Formula = "=D7/D5"

$wsWeapon.Range("A6").Value = "Exotic"
$wsWeapon.Range("B6").Value = 905
$wsWeapon.Range("C6").Value = 1000
$wsWeapon.Range("D6").Formula = "=(B6+C6)/2"
$wsWeapon.Range("E6").Formula = "=D7/D6"

$wsWeapon.Range("A7").Value = "Ascended"
$wsWeapon.Range("B7").Value = 950
$wsWeapon.Range("C7").Value = 1050
$wsWeapon.Range("D7").Formula = "=(B7+C7)/2"

# --- Second block (Coeff), columns F:J -------------------------------
$wsWeapon.Range("F1").Value = "Level"
$wsWeapon.Range("G1").Value = 35

$wsWeapon.Range("G2").Value = 253
$wsWeapon.Range("H2").Value = 279
$wsWeapon.Range("I2").Formula = "=(G2+H2)/2"
$wsWeapon.Range("J2").Formula = "=I5/I2"

$wsWeapon.Range("G4").Value = 344
$wsWeapon.Range("H4").Value = 380
$wsWeapon.Range("I4").Formula = "=(G4+H4)/2"
$wsWeapon.Range("J4").Formula = "=I5/I4"

$wsWeapon.Range("G5").Value = 373
$wsWeapon.Range("H5").Value = 412
$wsWeapon.Range("I5").Formula = "=(G5+H5)/2"

# --- Right summary block, columns L:O --------------------------------
$wsWeapon.Range("L1").Value = "Damage"
$wsWeapon.Range("M1").Value = 542
$wsWeapon.Range("O1").Value = 408

$wsWeapon.Range("L2").Value = "Coeff"
$wsWeapon.Range("M2").Value = 1.8
$wsWeapon.Range("O2").Value = 1.8

$wsWeapon.Range("L3").Value = "Power"
$wsWeapon.Range("M3").Value = 497
$wsWeapon.Range("O3").Value = 489

$wsWeapon.Range("L4").Value = "Armor"
$wsWeapon.Range("M4").Value = 624
$wsWeapon.Range("O4").Value = 624

$wsWeapon.Range("L5").Value = "Wep"
$wsWeapon.Range("M5").Formula = "=M1/M2/M3*M4"
$wsWeapon.Range("O5").Formula = "=O1/O2/O3*O4"

$wsWeapon.Range("L6").Value = "Expected"
$wsWeapon.Range("M6").Formula = "=G5*E5"
$wsWeapon.Range("O6").Formula = "=K5*G5"

$wsWeapon.Range("M7").Value = "Ascended"
$wsWeapon.Range("O7").Value = "Masterwork"

# --- Second, lower table (rows 10-16) --------------------------------
$wsWeapon.Range("A10").Value = "Level"
$wsWeapon.Range("B10").Value = 62
$wsWeapon.Range("L10").Value = "Damage"
$wsWeapon.Range("M10").Value = 1632

$wsWeapon.Range("B11").Value = 253
$wsWeapon.Range("C11").Value = 279
$wsWeapon.Range("D11").Formula = "=(B11+C11)/2"
$wsWeapon.Range("E11").Formula = "=D14/D11"
$wsWeapon.Range("L11").Value = "Coeff"
$wsWeapon.Range("M11").Value = 1.8

$wsWeapon.Range("L12").Value = "Power"
$wsWeapon.Range("M12").Value = 2482

$wsWeapon.Range("B13").Value = 597
$wsWeapon.Range("C13").Value = 659
$wsWeapon.Range("D13").Formula = "=(B13+C13)/2"
$wsWeapon.Range("E13").Formula = "=D14/D13"
$wsWeapon.Range("L13").Value = "Armor"
$wsWeapon.Range("M13").Value = 2597

$wsWeapon.Range("B14").Value = 643
$wsWeapon.Range("C14").Value = 710
$wsWeapon.Range("D14").Formula = "=(B14+C14)/2"
$wsWeapon.Range("E14").Formula = "=D15/D14"
$wsWeapon.Range("L14").Value = "Wep"
$wsWeapon.Range("M14").Formula = "=M10/M11/M12*M13"

$wsWeapon.Range("B15").Value = 727
$wsWeapon.Range("C15").Value = 803
$wsWeapon.Range("D15").Formula = "=(B15+C15)/2"
$wsWeapon.Range("L15").Value = "Expected"
$wsWeapon.Range("M15").Formula = "=I14*E14"

$wsWeapon.Range("M16").Value = "Ascended"

$wsWeapon.Range("M6").Select()
